# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match regenerated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 335
$ws1.Range("F4").Value  = 23
$ws1.Range("F5").Value  = 3248
$ws1.Range("F6").Value  = 2104
$ws1.Range("F7").Value  = 403
$ws1.Range("F8").Value  = 151
$ws1.Range("F9").Value  = 21
$ws1.Range("F10").Value = 1197
$ws1.Range("F11").Value = 216
$ws1.Range("F12").Value = 1157
$ws1.Range("F13").Value = 95

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 335
$ws4.Range("F4").Value  = 23
$ws4.Range("F5").Value  = 3248
$ws4.Range("F6").Value  = 2104
$ws4.Range("F7").Value  = 403
$ws4.Range("F9").Value  = 151
$ws4.Range("F10").Value = 21
$ws4.Range("F11").Value = 1197
$ws4.Range("F12").Value = 216
$ws4.Range("F13").Value = 1157
$ws4.Range("F14").Value = 95
